# "consistencia hasta la tabla de goles"
# "no me esta funcionando porque le falta id_rondas a la tabla de partidos"
#
# Edits the "tabla de triggers" on Sheet1:
#  - goles_x_jugador row (19): D19's rule text is corrected, and the row's
#    fill is made consistent with the rest of the table (entity-name cell
#    A19 reset to a plain white/background fill, rule cells C19/D19 given
#    the same yellow highlight already used on B19).
#  - planilla row (22): removes the stray "puede cambiar titular..." note
#    in C22 (no longer applicable).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Content fixes -------------------------------------------------------
$ws.Range("D19").Value = "no se puede fuera de fecha"
$ws.Range("C22").Value = ""

# --- Style consistency for row 19 ----------------------------------------
$ws.Range("C19").Interior.Color = $ws.Range("B19").Interior.Color
$ws.Range("D19").Interior.Color = $ws.Range("B19").Interior.Color
$ws.Range("A19").Interior.ThemeColor = 2
$ws.Range("A19").Interior.TintAndShade = 0

# --- Leave the selection where the author left it -------------------------
$ws.Range("C24").Select() | Out-Null
